$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 3
$ws.Range("B3").Value = "MCT-1A-Circuitos elétricos"
$ws.Range("C3").Value = "-"
$ws.Range("F3").Value = "MCT-3A-Máquinas Elétricas"

# Row 4
$ws.Range("B4").Value = "MCT-1A-Circuitos elétricos"
$ws.Range("C4").Value = "-"
$ws.Range("F4").Value = "MCT-3A-Máquinas Elétricas"

# Row 6
$ws.Range("B6").Value = "MCT-1A-Circuitos elétricos"
$ws.Range("C6").Value = "-"
$ws.Range("D6").Value = "-"

# Row 7
$ws.Range("B7").Value = "MCT-1A-Circuitos elétricos"
$ws.Range("C7").Value = "-"
$ws.Range("D7").Value = "-"
